$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "rm" variable row label first (matches author's original entry order,
# which places "rm" right after "variable" in the shared-strings table).
$ws.Range("A19").Value = "rm"

# New header row entries for the added columns.
$ws.Range("D1").Value = "p.value.w"
$ws.Range("E1").Value = "p.value.m"
$ws.Range("F1").Value = "p.value.m.eks"

# Correct the existing sys row (C8) p.value.sex value.
$ws.Range("C8").Value = 0.638

# Fill in the new D (p.value.w), E (p.value.m) and F (p.value.m.eks) columns.
$ws.Range("D2").Value = 0.649
$ws.Range("E2").Value = 0.564
$ws.Range("F2").Value = "NA"

$ws.Range("D3").Value = 0.578
$ws.Range("E3").Value = 0.677
$ws.Range("F3").Value = "NA"

$ws.Range("D4").Value = 0.22
$ws.Range("E4").Value = 0.089
$ws.Range("F4").Value = 0.423

$ws.Range("D5").Value = 0.263
$ws.Range("E5").Value = 0.0598
$ws.Range("F5").Value = 0.26

$ws.Range("D6").Value = 0.243
$ws.Range("E6").Value = 0.0714
$ws.Range("F6").Value = 0.322

$ws.Range("D7").Value = 0.292
$ws.Range("E7").Value = 0.248
$ws.Range("F7").Value = "NA"

$ws.Range("D8").Value = 0.583
$ws.Range("E8").Value = 0.231
$ws.Range("F8").Value = "NA"

$ws.Range("D9").Value = 0.0671
$ws.Range("E9").Value = 0.967
$ws.Range("F9").Value = "NA"

$ws.Range("D10").Value = 0.353
$ws.Range("E10").Value = 0.0306
$ws.Range("F10").Value = 0.12

$ws.Range("D11").Value = 0.781
$ws.Range("E11").Value = 0.68
$ws.Range("F11").Value = "NA"

$ws.Range("D12").Value = 0.219
$ws.Range("E12").Value = 0.396
$ws.Range("F12").Value = "NA"

$ws.Range("D13").Value = 0.324
$ws.Range("E13").Value = 0.0958
$ws.Range("F13").Value = "NA"

$ws.Range("D14").Value = 0.435
$ws.Range("E14").Value = 0.101
$ws.Range("F14").Value = 0.418

$ws.Range("D15").Value = 0.221
$ws.Range("E15").Value = 0.256
$ws.Range("F15").Value = "NA"

$ws.Range("D16").Value = 0.72
$ws.Range("E16").Value = 0.613
$ws.Range("F16").Value = "NA"

$ws.Range("D17").Value = 0.879
$ws.Range("E17").Value = 0.982
$ws.Range("F17").Value = "NA"

$ws.Range("D18").Value = 0.453
$ws.Range("E18").Value = 0.882
$ws.Range("F18").Value = "NA"

# New "rm" row (row 19) data.
$ws.Range("B19").Value = 0.243
$ws.Range("C19").Value = 0.588
$ws.Range("D19").Value = 0.614
$ws.Range("E19").Value = 0.396
$ws.Range("F19").Value = "NA"

# Column width adjustments: column A narrower, column F sized to fit its contents.
# (ColumnWidth is offset from the stored XML "width" by 5/6 of a character.)
$ws.Columns.Item(1).ColumnWidth = 11.998697916666666
$ws.Columns.Item(6).ColumnWidth = 11.498697916666666

# Leave the selection where the author left it when saving.
$ws.Range("E25").Select() | Out-Null
